$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows below the title (row 2), pushing everything else down by 3.
$ws.Rows("3:5").Insert()

# Row 3 stays blank (inherits style from the title row via Insert, same as before).
# Give it the same short height used by the author for the spacer row.
$ws.Rows(3).RowHeight = 15

# Row 4: short intro sentence about 2520.
$ws.Range("B4").Value = "2520 is the smallest number that can be divided by each of the numbers from 1 to 10 without any remainder."
$ws.Range("B4").Font.Name = "Aptos Narrow"
$ws.Rows(4).RowHeight = 15

# Row 5: the actual question, with "evenly divisible" emphasised (bold + underline).
$ws.Range("B5").Value = "What is the smallest positive number that is evenly divisible by all of the numbers from 1 to 20?"
$ws.Range("B5").Font.Name = "Aptos Narrow"
$emph = $ws.Range("B5").Characters(46, 16)
$emph.Font.Bold = $true
$emph.Font.Underline = $true
$ws.Rows(5).RowHeight = 15

# Row 6: blank spacer row before "1) Recursion".
$ws.Rows(6).RowHeight = 15

# Update the selection to match where the author ended up (question row).
$ws.Range("B19").Select()

Write-Output "edit applied"
